$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Citywide Totals")
$ws.Range("K2").Value = 2811
$ws.Range("K3").Value = 2731
$ws.Range("D4").Value = 1971
$ws.Range("K4").Value = 571
$ws.Range("K5").Value = 180
$ws.Range("K6").Value = 3357
$ws.Range("D7").Value = 28161
$ws.Range("K7").Value = 9650

$ws = $wb.Worksheets.Item("By Neighborhood")
$ws.Range("K2").Value = 74
$ws.Range("K6").Value = 80
$ws.Range("K7").Value = 284
$ws.Range("K8").Value = 634
$ws.Range("K10").Value = 53
$ws.Range("K15").Value = 98
$ws.Range("K18").Value = 66
$ws.Range("K19").Value = 284
$ws.Range("K20").Value = 220
$ws.Range("K22").Value = 31
$ws.Range("K29").Value = 501
$ws.Range("K33").Value = 378
$ws.Range("K34").Value = 47
$ws.Range("K36").Value = 110
$ws.Range("K37").Value = 315
$ws.Range("K40").Value = 22
$ws.Range("K42").Value = 339
$ws.Range("K43").Value = 84
$ws.Range("K44").Value = 91
$ws.Range("K48").Value = 115
$ws.Range("K52").Value = 269
$ws.Range("K53").Value = 139
$ws.Range("K54").Value = 182
$ws.Range("K55").Value = 104
$ws.Range("K59").Value = 17
$ws.Range("D63").Value = 350
$ws.Range("K63").Value = 38
$ws.Range("K65").Value = 227
$ws.Range("K67").Value = 377
$ws.Range("K69").Value = 23
$ws.Range("K71").Value = 31
$ws.Range("K74").Value = 13
$ws.Range("K76").Value = 145
$ws.Range("K77").Value = 69
$ws.Range("K79").Value = 247
$ws.Range("K80").Value = 32
$ws.Range("K83").Value = 210
$ws.Range("K84").Value = 68
$ws.Range("K85").Value = 463
$ws.Range("K88").Value = 110
$ws.Range("K89").Value = 128
$ws.Range("K92").Value = 39
$ws.Range("K94").Value = 115
$ws.Range("K95").Value = 157
$ws.Range("K96").Value = 134
$ws.Range("K98").Value = 55
$ws.Range("K99").Value = 173
$ws.Range("D101").Value = 28161
$ws.Range("K101").Value = 9650

$ws = $wb.Worksheets.Item("West Ridge")
$ws.Range("K3").Value = 21
$ws.Range("K7").Value = 134

$ws = $wb.Worksheets.Item("Auburn Gresham")
$ws.Range("K2").Value = 97
$ws.Range("K7").Value = 284

$ws = $wb.Worksheets.Item("Uptown")
$ws.Range("K3").Value = 41
$ws.Range("K4").Value = 19
$ws.Range("K7").Value = 128

$ws = $wb.Worksheets.Item("South Shore")
$ws.Range("K2").Value = 170
$ws.Range("K6").Value = 104
$ws.Range("K7").Value = 463

$ws = $wb.Worksheets.Item("Little Village")
$ws.Range("K2").Value = 73
$ws.Range("K3").Value = 64
$ws.Range("K6").Value = 112
$ws.Range("K7").Value = 269

$ws = $wb.Worksheets.Item("Norwood Park")
$ws.Range("K4").Value = 5
$ws.Range("K7").Value = 23

$ws = $wb.Worksheets.Item("Logan Square")
$ws.Range("K6").Value = 73
$ws.Range("K7").Value = 139

$ws = $wb.Worksheets.Item("Austin")
$ws.Range("K2").Value = 184
$ws.Range("K3").Value = 191
$ws.Range("K5").Value = 14
$ws.Range("K6").Value = 211
$ws.Range("K7").Value = 634

$ws = $wb.Worksheets.Item("South Chicago")
$ws.Range("K3").Value = 64
$ws.Range("K7").Value = 210

$ws = $wb.Worksheets.Item("Garfield Park")
$ws.Range("K2").Value = 105
$ws.Range("K3").Value = 137
$ws.Range("K4").Value = 20
$ws.Range("K7").Value = 378

$ws = $wb.Worksheets.Item("West Pullman")
$ws.Range("K3").Value = 54
$ws.Range("K7").Value = 157

$ws = $wb.Worksheets.Item("Grand Crossing")
$ws.Range("K2").Value = 82
$ws.Range("K3").Value = 107
$ws.Range("K7").Value = 315

$ws = $wb.Worksheets.Item("New City")
$ws.Range("K2").Value = 68
$ws.Range("K6").Value = 95
$ws.Range("K7").Value = 227

$ws = $wb.Worksheets.Item("Woodlawn")
$ws.Range("K4").Value = 10
$ws.Range("K6").Value = 44
$ws.Range("K7").Value = 173

$ws = $wb.Worksheets.Item("North Lawndale")
$ws.Range("K3").Value = 120
$ws.Range("K4").Value = 24
$ws.Range("K6").Value = 109
$ws.Range("K7").Value = 377

$ws = $wb.Worksheets.Item("South Deering")
$ws.Range("K3").Value = 23
$ws.Range("K7").Value = 68

$ws = $wb.Worksheets.Item("Loop")
$ws.Range("K2").Value = 35
$ws.Range("K3").Value = 57
$ws.Range("K7").Value = 182

$ws = $wb.Worksheets.Item("Englewood")
$ws.Range("K2").Value = 136
$ws.Range("K3").Value = 167
$ws.Range("K7").Value = 501

$ws = $wb.Worksheets.Item("Lake View")
$ws.Range("K2").Value = 19
$ws.Range("K7").Value = 115

$ws = $wb.Worksheets.Item("Chatham")
$ws.Range("K3").Value = 77
$ws.Range("K6").Value = 94
$ws.Range("K7").Value = 284

$ws = $wb.Worksheets.Item("Irving Park")
$ws.Range("K3").Value = 25
$ws.Range("K7").Value = 91

$ws = $wb.Worksheets.Item("River North")
$ws.Range("K3").Value = 25
$ws.Range("K6").Value = 86
$ws.Range("K7").Value = 145

$ws = $wb.Worksheets.Item("Ashburn")
$ws.Range("K2").Value = 28
$ws.Range("K7").Value = 80

$ws = $wb.Worksheets.Item("Humboldt Park")
$ws.Range("K2").Value = 86
$ws.Range("K3").Value = 109
$ws.Range("K7").Value = 339

$ws = $wb.Worksheets.Item("Avondale")
$ws.Range("K3").Value = 9
$ws.Range("K7").Value = 53

$ws = $wb.Worksheets.Item("Lower West Side")
$ws.Range("K3").Value = 25
$ws.Range("K6").Value = 36
$ws.Range("K7").Value = 104

$ws = $wb.Worksheets.Item("Roseland")
$ws.Range("K6").Value = 54
$ws.Range("K7").Value = 247

$ws = $wb.Worksheets.Item("Chicago Lawn")
$ws.Range("K2").Value = 77
$ws.Range("K3").Value = 59
$ws.Range("K7").Value = 220

$ws = $wb.Worksheets.Item("Calumet Heights")
$ws.Range("K2").Value = 20
$ws.Range("K4").Value = 9
$ws.Range("K7").Value = 66

$ws = $wb.Worksheets.Item("Grand Boulevard")
$ws.Range("K6").Value = 25
$ws.Range("K7").Value = 110

$ws = $wb.Worksheets.Item("Garfield Ridge")
$ws.Range("K3").Value = 15
$ws.Range("K7").Value = 47

$ws = $wb.Worksheets.Item("West Loop")
$ws.Range("K3").Value = 22
$ws.Range("K4").Value = 12
$ws.Range("K7").Value = 115

$ws = $wb.Worksheets.Item("Brighton Park")
$ws.Range("K3").Value = 23
$ws.Range("K7").Value = 98

$ws = $wb.Worksheets.Item("Wicker Park")
$ws.Range("K3").Value = 7
$ws.Range("K7").Value = 55

$ws = $wb.Worksheets.Item("Montclare")
$ws.Range("K3").Value = 4
$ws.Range("K7").Value = 17

$ws = $wb.Worksheets.Item("Albany Park")
$ws.Range("K6").Value = 27
$ws.Range("K7").Value = 74

$ws = $wb.Worksheets.Item("West Elsdon")
$ws.Range("K3").Value = 8
$ws.Range("K7").Value = 39

$ws = $wb.Worksheets.Item("United Center")
$ws.Range("K2").Value = 25
$ws.Range("K7").Value = 110

$ws = $wb.Worksheets.Item("Hyde Park")
$ws.Range("K3").Value = 23
$ws.Range("K6").Value = 36
$ws.Range("K7").Value = 84

$ws = $wb.Worksheets.Item("Clearing")
$ws.Range("K2").Value = 16
$ws.Range("K7").Value = 31

$ws = $wb.Worksheets.Item("Oakland")
$ws.Range("K2").Value = 14
$ws.Range("K7").Value = 31

$ws = $wb.Worksheets.Item("Riverdale")
$ws.Range("K2").Value = 34
$ws.Range("K3").Value = 24
$ws.Range("K7").Value = 69

$ws = $wb.Worksheets.Item("Rush & Division")
$ws.Range("K6").Value = 17
$ws.Range("K7").Value = 32

$ws = $wb.Worksheets.Item("Hegewisch")
$ws.Range("K3").Value = 11
$ws.Range("K7").Value = 22

$ws = $wb.Worksheets.Item("Printers Row")
$ws.Range("K3").Value = 3
$ws.Range("K7").Value = 13
